$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.598.27"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "1.850.27"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'314.19"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  -2.40%  "
$ws.Range("D8").Value = "'0.3642"
$ws.Range("E8").Value = "  -3.10%  "
$ws.Range("D9").Value = "'44.69"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").Value = "'0.07306"
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("D11").Value = "'0.8763"
$ws.Range("E11").Value = "  -6.51%  "
$ws.Range("D12").Value = "'20.73"
$ws.Range("E12").Value = "  -2.60%  "
$ws.Range("D13").Value = "1.907.38"
$ws.Range("E13").Value = "  +3.81%  "
$ws.Range("D14").Value = "'5.347"
$ws.Range("E14").Value = "  -1.80%  "
$ws.Range("D15").Value = "'6.535"
$ws.Range("E15").Value = "  -3.11%  "
$ws.Range("D16").Value = "'0.06923"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "'1.005"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "'78.88"
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("D19").Value = "'0.000008882"
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "'15.39"
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("D22").Value = "27.618.70"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").Value = "'5.004"
$ws.Range("E23").Value = "  -2.36%  "
$ws.Range("D24").Value = "'10.63"
$ws.Range("E24").Value = "  -3.86%  "
$ws.Range("D25").Value = "2.144.78"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").Value = "'1.984"
$ws.Range("E26").Value = "  -2.57%  "
$ws.Range("D27").Value = "'153.47"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").Value = "'19.07"
$ws.Range("E28").Value = "  +2.36%  "
$ws.Range("D29").Value = "'121.15"
$ws.Range("E29").Value = "  +6.72%  "
$ws.Range("D30").Value = "'5.266"
$ws.Range("E30").Value = "  -6.00%  "
$ws.Range("E31").Value = "  +11.90%  "
$ws.Range("D32").Value = "'0.08925"
$ws.Range("D33").Value = "'0.7616"
$ws.Range("E33").Value = "  -6.23%  "
$ws.Range("D34").Value = "'4.575"
$ws.Range("E34").Value = "  -4.86%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'1.100"
$ws.Range("E36").Value = "  -7.12%  "
$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "'0.05373"
$ws.Range("E38").Value = "  -2.66%  "
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("E41").Value = "  -4.92%  "
$ws.Range("D42").Value = "'6.924"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").Value = "'0.5106"
$ws.Range("E43").Value = "  -3.14%  "
$ws.Range("D44").Value = "'0.1650"
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("D45").Value = "'8.273"
$ws.Range("E45").Value = "  -5.76%  "
$ws.Range("D46").Value = "'0.06562"
$ws.Range("E46").Value = "  -2.87%  "
$ws.Range("D47").Value = "'10.41"
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("D48").Value = "'0.4740"
$ws.Range("E48").Value = "  -3.11%  "
$ws.Range("D49").Value = "'104.45"
$ws.Range("E49").Value = "  -2.28%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").Value = "'1.627"
$ws.Range("E51").Value = "  -2.93%  "
